# Update crypto price/volume snapshot values (GitHub Actions daily refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BNB
$ws.Range("E2").Value = "'0.85%"

# Row 3: OKB
$ws.Range("D3").Value = "'27.03"
$ws.Range("E3").Value = "'0.72%"

# Row 5: Cronos
$ws.Range("D5").Value = "'0.06181"
$ws.Range("E5").Value = "'3.28%"

# Row 6: KuCoinToken
$ws.Range("D6").Value = "'6.686"
$ws.Range("E6").Value = "'0.69%"

# Row 7: MXToken
$ws.Range("D7").Value = "'0.8500"
$ws.Range("E7").Value = "'-0.88%"

# Row 8: FTXToken
$ws.Range("D8").Value = "'0.9169"
$ws.Range("E8").Value = "'-0.50%"

# Row 9: WazirX
$ws.Range("D9").Value = "'0.1405"
$ws.Range("E9").Value = "'1.24%"

# Row 10: LiechtensteinCryptoassetsExchange
$ws.Range("D10").Value = "'0.04672"
$ws.Range("E10").Value = "'4.23%"

# Row 11: MandalaExchangeToken
$ws.Range("D11").Value = "'0.07087"
$ws.Range("E11").Value = "'0.95%"

# Row 12: BitrueCoin
$ws.Range("D12").Value = "'0.03157"
$ws.Range("E12").Value = "'3.45%"

# Row 13: BitMartToken
$ws.Range("E13").Value = "'-0.54%"

# Row 14: BitForexToken
$ws.Range("D14").Value = "'0.001525"
$ws.Range("E14").Value = "'-0.97%"

# Row 15: One
$ws.Range("D15").Value = "'0.0006176"
$ws.Range("E15").Value = "'2.12%"

# Row 16: TigerCash
$ws.Range("D16").Value = "'0.006117"
$ws.Range("E16").Value = "'-0.29%"

# Row 17: LEO
$ws.Range("D17").Value = "'3.451"
$ws.Range("E17").Value = "'0.54%"

# Row 18: GateToken
$ws.Range("D18").Value = "'3.177"
$ws.Range("E18").Value = "'0.90%"

# Row 19: BTSEToken
$ws.Range("D19").Value = "'2.168"
$ws.Range("E19").Value = "'0.69%"

# Row 21: ProBitToken
$ws.Range("D21").Value = "'0.1299"
$ws.Range("E21").Value = "'0.86%"

# Row 22: MCDex
$ws.Range("D22").Value = "'4.081"
$ws.Range("E22").Value = "'1.41%"

# Row 23: CoinExToken
$ws.Range("D23").Value = "'0.04234"
$ws.Range("E23").Value = "'0.18%"

# Row 24: BitKan
$ws.Range("D24").Value = "'0.001217"
$ws.Range("E24").Value = "'-0.02%"

# Row 25: HotbitToken
$ws.Range("E25").Value = "'-5.64%"

# Row 26: NitroEx
$ws.Range("E26").Value = "'0.10%"

# Row 40: IDEX
$ws.Range("D40").Value = "'0.03885"
$ws.Range("E40").Value = "'1.49%"

# Row 41: BKEXToken
$ws.Range("E41").Value = "'-0.10%"

# Row 42: KickToken
$ws.Range("D42").Value = "'0.004100"
$ws.Range("E42").Value = "'9.08%"

# Row 43: LocalTraders
$ws.Range("E43").Value = "'7.34%"

# Row 44: CEJI
$ws.Range("E44").Value = "'-10.06%"

# Row 45: CoinLion
$ws.Range("D45").Value = "'0.00005162"
$ws.Range("E45").Value = "'0.65%"

# Row 46: Kangarootoken
$ws.Range("E46").Value = "'0.08%"

# Row 48: BOLO
$ws.Range("D48").Value = "'0.1667"
$ws.Range("E48").Value = "'52.61%"

# Row 49: CryptobidCoin
$ws.Range("E49").Value = "'0.08%"

# Row 50: SpecialPowerGold
$ws.Range("E50").Value = "'0.08%"
